$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.142.50"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.675.95"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.35"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.85"
$ws.Range("E8").Value = "  +6.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.262"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0621"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "1.911.57"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "1.672.00"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.21"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.560"
$ws.Range("E15").Value = "  +4.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.50"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "27.100.84"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.15"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").Value = "0.0₃0742"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.80"
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.54"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.56"
$ws.Range("E23").Value = "  +3.13%  "
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.06"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.46"
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.41"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("D33").Value = "1.543.48"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.66"
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.946"
$ws.Range("E37").Value = "  +3.27%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.99"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("E42").Value = "  +4.33%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "1.820.61"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.64"
$ws.Range("E47").Value = "  +6.47%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.66"
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.19"
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("E51").Value = "  -0.12%  "
